# Apply the Fri Nov  1 13:55:51 UTC 2024 cryptos-list refresh (GitHub Actions scrape update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: target cell reference + the new text that belongs there.
# "NumericLooking" cells (plain decimals like "0.999") get their number
# format forced to Text first so Excel does not silently coerce the
# assigned string into a Double (which would mangle values such as
# "1.00" -> 1 or introduce floating point noise). The format is then
# restored to the sheet's normal (General) style so no stray cell
# formatting is left behind.
$updates = @(
    @{ Cell = 'D2'; Value = '70.899.34'; NumericLooking = $false }
    @{ Cell = 'E2'; Value = '  -0.83%  '; NumericLooking = $false }
    @{ Cell = 'D3'; Value = '2.561.90'; NumericLooking = $false }
    @{ Cell = 'E3'; Value = '  -1.86%  '; NumericLooking = $false }
    @{ Cell = 'D4'; Value = '0.999'; NumericLooking = $true }
    @{ Cell = 'E4'; Value = '  -0.14%  '; NumericLooking = $false }
    @{ Cell = 'D5'; Value = '583.88'; NumericLooking = $true }
    @{ Cell = 'E5'; Value = '  -0.09%  '; NumericLooking = $false }
    @{ Cell = 'D6'; Value = '172.21'; NumericLooking = $true }
    @{ Cell = 'E6'; Value = '  -0.21%  '; NumericLooking = $false }
    @{ Cell = 'E7'; Value = '  +0.02%  '; NumericLooking = $false }
    @{ Cell = 'E8'; Value = '  +0.86%  '; NumericLooking = $false }
    @{ Cell = 'D9'; Value = '2.558.29'; NumericLooking = $false }
    @{ Cell = 'E9'; Value = '  -2.03%  '; NumericLooking = $false }
    @{ Cell = 'D10'; Value = '0.166'; NumericLooking = $true }
    @{ Cell = 'E10'; Value = '  -0.63%  '; NumericLooking = $false }
    @{ Cell = 'E11'; Value = '  -0.94%  '; NumericLooking = $false }
    @{ Cell = 'D12'; Value = '0.357'; NumericLooking = $true }
    @{ Cell = 'E12'; Value = '  +0.21%  '; NumericLooking = $false }
    @{ Cell = 'D13'; Value = '4.93'; NumericLooking = $true }
    @{ Cell = 'E13'; Value = '  +0.81%  '; NumericLooking = $false }
    @{ Cell = 'D14'; Value = '3.026.08'; NumericLooking = $false }
    @{ Cell = 'E14'; Value = '  -2.07%  '; NumericLooking = $false }
    @{ Cell = 'D15'; Value = '70.591.39'; NumericLooking = $false }
    @{ Cell = 'E15'; Value = '  -1.27%  '; NumericLooking = $false }
    @{ Cell = 'D16'; Value = '0.0000179'; NumericLooking = $true }
    @{ Cell = 'E16'; Value = '  -4.14%  '; NumericLooking = $false }
    @{ Cell = 'D17'; Value = '25.55'; NumericLooking = $true }
    @{ Cell = 'E17'; Value = '  +0.05%  '; NumericLooking = $false }
    @{ Cell = 'D18'; Value = '2.549.85'; NumericLooking = $false }
    @{ Cell = 'E18'; Value = '  -3.51%  '; NumericLooking = $false }
    @{ Cell = 'D19'; Value = '7.94'; NumericLooking = $true }
    @{ Cell = 'E19'; Value = '  +1.63%  '; NumericLooking = $false }
    @{ Cell = 'D20'; Value = '11.49'; NumericLooking = $true }
    @{ Cell = 'E20'; Value = '  -4.60%  '; NumericLooking = $false }
    @{ Cell = 'D21'; Value = '357.35'; NumericLooking = $true }
    @{ Cell = 'E21'; Value = '  -3.60%  '; NumericLooking = $false }
    @{ Cell = 'D22'; Value = '3.95'; NumericLooking = $true }
    @{ Cell = 'E22'; Value = '  -2.45%  '; NumericLooking = $false }
    @{ Cell = 'D23'; Value = '2.07'; NumericLooking = $true }
    @{ Cell = 'E23'; Value = '  +2.83%  '; NumericLooking = $false }
    @{ Cell = 'D24'; Value = '1.00'; NumericLooking = $true }
    @{ Cell = 'E24'; Value = '  +0.04%  '; NumericLooking = $false }
    @{ Cell = 'D25'; Value = '70.37'; NumericLooking = $true }
    @{ Cell = 'E25'; Value = '  -1.41%  '; NumericLooking = $false }
    @{ Cell = 'D26'; Value = '4.08'; NumericLooking = $true }
    @{ Cell = 'E26'; Value = '  -2.41%  '; NumericLooking = $false }
    @{ Cell = 'D27'; Value = '9.20'; NumericLooking = $true }
    @{ Cell = 'E27'; Value = '  -1.32%  '; NumericLooking = $false }
    @{ Cell = 'D29'; Value = '0.993'; NumericLooking = $true }
    @{ Cell = 'E29'; Value = '  -0.82%  '; NumericLooking = $false }
    @{ Cell = 'D30'; Value = '0.0₃0928'; NumericLooking = $false }
    @{ Cell = 'E30'; Value = '  -1.33%  '; NumericLooking = $false }
    @{ Cell = 'D31'; Value = '7.97'; NumericLooking = $true }
    @{ Cell = 'E31'; Value = '  +0.57%  '; NumericLooking = $false }
    @{ Cell = 'B32'; Value = 'Bittensor'; NumericLooking = $false }
    @{ Cell = 'C32'; Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'; NumericLooking = $false }
    @{ Cell = 'D32'; Value = '475.19'; NumericLooking = $true }
    @{ Cell = 'E32'; Value = '  -1.99%  '; NumericLooking = $false }
    @{ Cell = 'B33'; Value = 'Fetch.AI'; NumericLooking = $false }
    @{ Cell = 'C33'; Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'; NumericLooking = $false }
    @{ Cell = 'D33'; Value = '1.29'; NumericLooking = $true }
    @{ Cell = 'E33'; Value = '  -1.80%  '; NumericLooking = $false }
    @{ Cell = 'E34'; Value = '  -0.53%  '; NumericLooking = $false }
    @{ Cell = 'E35'; Value = '  -0.05%  '; NumericLooking = $false }
    @{ Cell = 'D36'; Value = '0.120'; NumericLooking = $true }
    @{ Cell = 'E36'; Value = '  +4.27%  '; NumericLooking = $false }
    @{ Cell = 'D37'; Value = '159.58'; NumericLooking = $true }
    @{ Cell = 'E37'; Value = '  +0.08%  '; NumericLooking = $false }
    @{ Cell = 'D38'; Value = '19.10'; NumericLooking = $true }
    @{ Cell = 'E38'; Value = '  +1.03%  '; NumericLooking = $false }
    @{ Cell = 'D39'; Value = '18.76'; NumericLooking = $true }
    @{ Cell = 'E39'; Value = '  -2.31%  '; NumericLooking = $false }
    @{ Cell = 'E40'; Value = '  +0.03%  '; NumericLooking = $false }
    @{ Cell = 'D41'; Value = '4.90'; NumericLooking = $true }
    @{ Cell = 'E41'; Value = '  +1.58%  '; NumericLooking = $false }
    @{ Cell = 'E42'; Value = '  -0.18%  '; NumericLooking = $false }
    @{ Cell = 'D43'; Value = '1.63'; NumericLooking = $true }
    @{ Cell = 'E43'; Value = '  -4.69%  '; NumericLooking = $false }
    @{ Cell = 'D44'; Value = '2.38'; NumericLooking = $true }
    @{ Cell = 'E44'; Value = '  -6.31%  '; NumericLooking = $false }
    @{ Cell = 'E45'; Value = '  -12.94%  '; NumericLooking = $false }
    @{ Cell = 'D46'; Value = '38.49'; NumericLooking = $true }
    @{ Cell = 'E46'; Value = '  -0.93%  '; NumericLooking = $false }
    @{ Cell = 'D47'; Value = '145.81'; NumericLooking = $true }
    @{ Cell = 'E47'; Value = '  -2.00%  '; NumericLooking = $false }
    @{ Cell = 'D48'; Value = '0.541'; NumericLooking = $true }
    @{ Cell = 'E48'; Value = '  +0.59%  '; NumericLooking = $false }
    @{ Cell = 'D49'; Value = '3.56'; NumericLooking = $true }
    @{ Cell = 'E49'; Value = '  -1.91%  '; NumericLooking = $false }
    @{ Cell = 'E50'; Value = '  -1.70%  '; NumericLooking = $false }
    @{ Cell = 'D51'; Value = '0.0742'; NumericLooking = $true }
    @{ Cell = 'E51'; Value = '  +0.10%  '; NumericLooking = $false }
)

foreach ($update in $updates) {
    $range = $ws.Range($update.Cell)
    if ($update.NumericLooking) {
        # Force text storage so "1.00"/"0.999"-style values keep their
        # exact textual form instead of becoming numbers.
        $range.NumberFormat = "@"
        $range.Value = $update.Value
        $range.Style = "Normal"
    } else {
        $range.Value = $update.Value
    }
}
